$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.588.22"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.476.69"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'0.9563"
$ws.Range("E5").Value = "  +6.85%  "
$ws.Range("D6").Value = "'280.16"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.3653"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("D8").Value = "'0.3064"
$ws.Range("E8").Value = "  -4.38%  "
$ws.Range("D9").Value = "'39.96"
$ws.Range("E9").Value = "  -1.08%  "
$ws.Range("D10").Value = "'1.061"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").Value = "'0.06674"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'5.520"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'18.06"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "'6.218"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.00001035"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "'0.9577"
$ws.Range("E17").Value = "  +6.62%  "
$ws.Range("D18").Value = "1.476.85"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'0.05945"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").Value = "'69.96"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'5.500"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("D22").Value = "'14.44"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").Value = "'11.06"
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").Value = "'2.262"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("D25").Value = "20.626.27"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'143.21"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").Value = "'2.113"
$ws.Range("E27").Value = "  -8.09%  "
$ws.Range("D28").Value = "'17.27"
$ws.Range("D29").Value = "1.637.67"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'113.87"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'3.965"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'5.012"
$ws.Range("E32").Value = "  -4.80%  "
$ws.Range("D33").Value = "'0.8137"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("D34").Value = "'0.07948"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").Value = "'1.513"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").Value = "'1.218"
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("D37").Value = "'0.05821"
$ws.Range("E37").Value = "  -5.70%  "
$ws.Range("D38").Value = "'4.738"
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").Value = "'0.02049"
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").Value = "'0.9578"
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("D41").Value = "'10.38"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").Value = "'0.1878"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").Value = "'7.454"
$ws.Range("E43").Value = "  +3.81%  "
$ws.Range("D44").Value = "'0.5305"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").Value = "'3.542"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").Value = "'12.23"
$ws.Range("E46").Value = "  -1.77%  "
$ws.Range("D47").Value = "'118.01"
$ws.Range("E47").Value = "  -4.49%  "
$ws.Range("D48").Value = "'0.5188"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("D49").Value = "'1.822"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "'0.06487"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("D51").Value = "'0.9868"
$ws.Range("E51").Value = "  -0.57%  "
